# Commit: "Finally removed config sync and family, fixed some bugs relating
# to the kit cache and other minor kit bugs."
#
# Data change: the netcalls table lists the ConfigSync.NetCalls and
# ListSync.NetCalls families (IDs 3000-3005, rows 101-106) which were
# removed from the codebase. Delete those six table rows; Excel's table
# (Table2) auto-shrinks from A1:E116 to A1:E110, and every row below
# shifts up to fill the gap (old rows 107-116 become new rows 101-110).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sanity check before deleting - row 101 should be the first of the
# ListSync.NetCalls / ConfigSync.NetCalls block (ID 3000) and row 106
# should be the last of that block (ID 3005) before the KitSync.NetCalls
# block (ID 3006+) continues on what will become the new row 101.
$firstId = $ws.Range("B101").Value2
$lastId  = $ws.Range("B106").Value2

if ($firstId -eq 3000 -and $lastId -eq 3005) {
    $ws.Range("A101:E106").EntireRow.Delete() | Out-Null
} else {
    throw "Expected rows 101-106 to be the ConfigSync/ListSync NetCalls block (IDs 3000-3005), but found IDs $firstId..$lastId. Aborting to avoid deleting the wrong rows."
}

# Move the selection to where the author's saved view left it.
$ws.Range("A101").Select() | Out-Null
